$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(6, 1).Value = 5

    # "2026-02-16" looks like a date, so Excel would silently convert it to a
    # date serial number on assignment. Force the cell to Text first so the
    # literal string is preserved, then strip the formatting back off so the
    # cell ends up with the default (unstyled) look, matching the other rows.
    $ws.Cells.Item(6, 2).NumberFormat = "@"
    $ws.Cells.Item(6, 2).Value = "2026-02-16"
    $ws.Cells.Item(6, 2).ClearFormats()

    $ws.Cells.Item(6, 3).Value = "22:56:50"
    $ws.Cells.Item(6, 4).Value = "base_strategy"
    $ws.Cells.Item(6, 5).Value = "DOWN"
    $ws.Cells.Item(6, 6).Value = 0.5

    # Exit Price stays blank (trade still open) but the cell itself must
    # exist in the sheet, so nudge it with a no-op formatting call.
    $ws.Cells.Item(6, 7).Borders.LineStyle = -4142

    $ws.Cells.Item(6, 8).Value = "OPEN"
    $ws.Cells.Item(6, 9).Value = 0
    $ws.Cells.Item(6, 10).Value = 0
    $ws.Cells.Item(6, 11).Value = 100
    $ws.Cells.Item(6, 12).Value = 0
    $ws.Cells.Item(6, 13).Value = 0
    $ws.Cells.Item(6, 14).Value = 0.6
    $ws.Cells.Item(6, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason also stays blank while the trade is open.
    $ws.Cells.Item(6, 16).Borders.LineStyle = -4142

    $ws.Cells.Item(6, 17).Value = 0
}
